$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value. Column D numeric-looking values are
# given a leading apostrophe so Excel keeps them as text (matching the
# original inlineStr cells) instead of silently converting to Double.
$updates = [ordered]@{
    'D2' = '61.759.10'
    'E2' = '  +1.89%  '
    'D3' = '2.412.32'
    'E3' = '  +3.77%  '
    'E4' = '  -1.07%  '
    'D5' = "'553.27"
    'E5' = '  +2.21%  '
    'D6' = "'142.26"
    'E6' = '  +5.22%  '
    'D7' = "'0.999"
    'E7' = '  +0.66%  '
    'D8' = "'0.523"
    'E8' = '  +0.15%  '
    'D9' = '2.416.05'
    'E9' = '  +3.29%  '
    'E10' = '  +4.03%  '
    'E11' = '  +1.33%  '
    'D12' = "'5.37"
    'E12' = '  +1.87%  '
    'D13' = "'0.352"
    'E13' = '  +3.98%  '
    'D14' = "'26.08"
    'E14' = '  +7.04%  '
    'D15' = "'0.0000175"
    'E15' = '  +10.08%  '
    'D16' = '2.847.82'
    'E16' = '  +3.33%  '
    'D17' = '61.526.64'
    'E17' = '  +24.45%  '
    'D18' = '2.412.56'
    'E18' = '  +8.39%  '
    'D19' = "'11.14"
    'E19' = '  +5.92%  '
    'D20' = "'323.29"
    'E20' = '  +3.16%  '
    'D21' = "'4.17"
    'E21' = '  +2.63%  '
    'D22' = "'6.71"
    'E22' = '  +2.99%  '
    'E23' = '  +0.41%  '
    'D24' = "'64.46"
    'E24' = '  +3.02%  '
    'D25' = "'1.75"
    'E25' = '  +3.51%  '
    'D26' = "'9.18"
    'E26' = '  +10.50%  '
    'D27' = "'559.41"
    'E27' = '  +12.24%  '
    'E28' = '  +0.14%  '
    'D29' = '2.500.76'
    'D30' = "'8.31"
    'E30' = '  +5.90%  '
    'D31' = '0.0₃0921'
    'E31' = '  +6.60%  '
    'D32' = "'1.44"
    'E32' = '  +5.39%  '
    'E34' = '  +4.29%  '
    'D35' = "'1.54"
    'E35' = '  +2.85%  '
    'D36' = "'5.79"
    'E36' = '  +12.04%  '
    'E37' = '  +0.25%  '
    'E38' = '  +11.92%  '
    'D39' = "'4.77"
    'E39' = '  +4.41%  '
    'D40' = "'0.384"
    'E40' = '  +3.77%  '
    'D41' = "'18.70"
    'E41' = '  +1.61%  '
    'D42' = "'146.51"
    'E42' = '  +4.24%  '
    'D43' = "'0.999"
    'E43' = '  +0.01%  '
    'B44' = 'Aave'
    'C44' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D44' = "'148.98"
    'E44' = '  +7.04%  '
    'B45' = 'dogwifhat'
    'C45' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D45' = "'2.25"
    'E45' = '  +9.79%  '
    'E46' = '  +3.12%  '
    'D47' = "'0.0535"
    'E47' = '  +5.32%  '
    'D48' = "'20.26"
    'E48' = '  +6.27%  '
    'D49' = "'0.591"
    'E49' = '  +4.64%  '
    'D50' = "'0.0909"
    'E50' = '  +2.04%  '
    'D51' = "'0.0225"
    'E51' = '  +2.77%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# Cells that received the leading apostrophe got an implicit "quote prefix"
# text style applied; reset them back to the workbook default "Normal" style
# so no stray cell formatting is introduced.
foreach ($ref in $updates.Keys) {
    if ($updates[$ref].ToString().StartsWith("'")) {
        $ws.Range($ref).Style = "Normal"
    }
}
